# Update ranking at 2025-12-04 09:13
# Append a new row (row 49) to the bottom of the tracking table with a new
# timestamp entry and placeholder "-" values for the other columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 49

$ws.Cells.Item($newRow, 1).Value = "2025/12/04 18:00"
$ws.Cells.Item($newRow, 2).Value = "-"
$ws.Cells.Item($newRow, 3).Value = "-"
$ws.Cells.Item($newRow, 4).Value = "-"
$ws.Cells.Item($newRow, 5).Value = "-"
$ws.Cells.Item($newRow, 6).Value = "-"
$ws.Cells.Item($newRow, 7).Value = "-"
